# Update master to output generated at c986bee
$d = $word.ActiveDocument

$replacements = @(
    @("2024-12-04 Wednesday", "2024-12-05 Thursday"),
    @("871×8=", "718×5="),
    @("576×6=", "519×7="),
    @("296×9=", "618×5="),
    @("925×6=", "137×5="),
    @("863×5=", "307×7="),
    @("348×2=", "676×5="),
    @("442×6=", "689×3="),
    @("242×5=", "897×8="),
    @("592×7=", "157×5="),
    @("990×6=", "901×8="),
    @("277×2=", "792×9="),
    @("438×2=", "809×7="),
    @("781×7=", "286×7="),
    @("838×5=", "495×6="),
    @("370×8=", "605×6="),
    @("405×9=", "858×2="),
    @("454×6=", "232×7="),
    @("417×9=", "455×8="),
    @("847×3=", "946×8="),
    @("633×4=", "342×2="),
    @("300×5=", "251×6="),
    @("779×5=", "679×2="),
    @("874×8=", "506×8="),
    @("152×5=", "684×3="),
    @("849×6=", "218×3=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

$d.Save()
